$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (was date 44204 / Primera-110) becomes date 44189 / Especial-20
$ws.Range("D2").Value = 44189
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 2143

# Row 3 (date stays 44189) becomes Primera-30
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 13000
$ws.Range("S3").Value = 1857

# Row 4 (was date 44189 / Primera-30) becomes date 44204 / Primera-110
$ws.Range("D4").Value = 44204
$ws.Range("M4").Value = 110
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7500
$ws.Range("P4").Value = 7318
$ws.Range("S4").Value = 1045
